$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (tomlora) updated stats
$ws.Range("D6").Value = 117
$ws.Range("E6").Value = 104
$ws.Range("F6").Value = 58.42500000000002
$ws.Range("G6").Value = 2401
$ws.Range("H6").Value = 1014
$ws.Range("I6").Value = 335
$ws.Range("J6").Value = 322
$ws.Range("K6").Value = 23993
$ws.Range("L6").Value = 988
$ws.Range("M6").Value = 718
$ws.Range("N6").Value = 895
$ws.Range("O6").Value = 8.444444444444445
$ws.Range("P6").Value = 6.136752136752137
$ws.Range("Q6").Value = 7.64957264957265
$ws.Range("R6").Value = 20.52
$ws.Range("S6").Value = 29.96

# Row 10 (nukethestars) updated stats
$ws.Range("D10").Value = 22
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = 11.26316666666667
$ws.Range("G10").Value = 1647
$ws.Range("H10").Value = 678
$ws.Range("I10").Value = 239
$ws.Range("J10").Value = 166
$ws.Range("K10").Value = 900
$ws.Range("L10").Value = 70
$ws.Range("M10").Value = 163
$ws.Range("N10").Value = 293
$ws.Range("O10").Value = 3.181818181818182
$ws.Range("P10").Value = 7.409090909090909
$ws.Range("Q10").Value = 13.31818181818182
$ws.Range("R10").Value = 74.86
$ws.Range("S10").Value = 30.72
